# Settings back to qfin22 presentation - recalibration
# Updates the calibration parameter values (and their autocorrelation
# arrays) on the AR, SETAR, GARCH, TARCH and AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# ---- AR sheet ----
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.02613277658711496
$ws.Range("B3").Value = 0.8301829873809643
$ws.Range("B4").Value = 32.32452110185174
$ws.Range("B5").Value = "[1.0, 0.275705054056148, 0.2208689431034171, 0.2446441446296424, 0.21965573342552783, 0.45039648787823416, 0.22322631753267497, 0.17829456991036335, 0.20361838270477414, 0.16864611601188503, 0.16274842939876874, 0.11244403587872488, 0.12681062917433916, 0.16820412158963902, 0.1725092612066047, 0.13052782616264305, 0.12997100497246106, 0.12184420457526338, 0.13471377414492197, 0.13825366197893893]"

# ---- SETAR sheet ----
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -1.648913426514714
$ws.Range("B5").Value = 0.7691280218622456
$ws.Range("B6").Value = 25.88846556954054
$ws.Range("B7").Value = 1.467849337606879
$ws.Range("B8").Value = 0.7692302244151693
$ws.Range("B9").Value = 21.98748635576831
$ws.Range("B10").Value = "[1.0, 0.2739530517359736, 0.2224915077367662, 0.2536163616061678, 0.24307049246815776, 0.3371226102159966, 0.2013129362220637, 0.16768742913364268, 0.17404639480184256, 0.16401359531854975, 0.16960940966295726, 0.12399317532924053, 0.12755655395776289, 0.1341383400555717, 0.13251755700719617, 0.11120849860614665, 0.1308505205751499, 0.08973544661224055, 0.11059761629985555, 0.14218000607443212]"

# ---- GARCH sheet ----
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.004578046838512801
$ws.Range("B3").Value = 0.2711836398779733
$ws.Range("B4").Value = 0.1154968220522724
$ws.Range("B5").Value = 0.8845031787451064
$ws.Range("B6").Value = "[1.0, 0.017954617471897607, -0.02269678883968538, -0.016664094552898457, -0.03995012981988733, 0.30743984014274184, 0.002661443544310895, 0.0008193156949588996, 0.0068745078158507795, -0.040667503721164365, -0.06303825613328595, -0.06975543440111025, -0.04480606904798816, -0.003656602177575346, 0.011974541158367274, -0.027434005015989013, -0.032327166946387646, -0.04714387313833362, -0.043138765816501695, -0.04973842852195918]"

# ---- TARCH sheet ----
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.04900101981963107
$ws.Range("B3").Value = 0.2796132217125263
$ws.Range("B4").Value = 0.09590111486823154
$ws.Range("B5").Value = 0.8837189556014422
$ws.Range("B6").Value = "[1.0, 0.016293295128490586, -0.022296880341892932, -0.016771700679996317, -0.04070812889552454, 0.3068440952239229, 0.003985903777195831, 0.002028727930121994, 0.007509245779219818, -0.040328067073641914, -0.06352310480717599, -0.06910382110981665, -0.044537964711827895, -0.002720803709117024, 0.01377814692061423, -0.02603552927779307, -0.03196739560762685, -0.047133636081276006, -0.0429014761511698, -0.04858436592152345]"
$ws.Range("B7").Value = 0.0407598595947847

# ---- AR_TARCH sheet ----
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.06427563470148839
$ws.Range("B3").Value = 0.2586096398149271
$ws.Range("B4").Value = 0.1102726332407875
$ws.Range("B5").Value = 0.8860174228647901
$ws.Range("B6").Value = "[1.0, 0.02845062759557058, -0.013212133038455878, 0.012091315094959683, -0.03250494283638258, 0.2469566830443298, 0.006053571416459681, -0.0019865960400446015, -0.005487518425938155, -0.04645540813986593, -0.035317120333900796, -0.05756204350142196, -0.03611024729303446, -0.015211551946075798, 0.011112813821146508, -0.040169679848970026, -0.03295441395979722, -0.037048078646522854, -0.03312443051899934, -0.034632428944231214]"
$ws.Range("B7").Value = 0.007419887593253347
$ws.Range("B9").Value = 0.8077329722011638
